# Add a new service-event row (row 20) to the "Card13" sheet.
# Commit: "إضافة حدث جديد في Card13 by admin at 2025-12-08 11:41:41"
#   -> A new event-log entry is appended right after the last existing
#      event (row 19), extending the sheet's used range from A1:O19 to
#      A1:O20.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card13")

$row = 20

# Same "card" id as every other row on this sheet - stored as text like
# the rest of column A.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "13"

# Columns B-K (tonnage range / checklist marks) are left blank for this
# event, matching every other pure "event" row (14-19) on the sheet.

# Date of the service event (column L).
$ws.Cells.Item($row, 12).Value = "24\7\2025"

# Event / correction description (column M = "Correction").
$ws.Cells.Item($row, 13).Value = "تم تغيير الفلاتس وتغيير جريد 1"

# Column N ("Event") is left blank for this entry.

# Serviced by (column O).
$ws.Cells.Item($row, 15).Value = "الخبير"
